# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos - updated totals
$ws.Range("B4").Value = 668174
$ws.Range("C4").Value = 20026
$ws.Range("D4").Value = 57228
$ws.Range("E4").Value = 577015
$ws.Range("G4").Value = 1343
$ws.Range("H4").Value = 33931

# Row 7: Francia - updated totals
$ws.Range("B7").Value = 165027
$ws.Range("C7").Value = 17164
$ws.Range("E7").Value = 114295

# Row 15: Canada - updated totals
$ws.Range("B15").Value = 29929
$ws.Range("C15").Value = 1550
$ws.Range("D15").Value = 9674
$ws.Range("E15").Value = 19064

# Row 20: Austria - updated totals
$ws.Range("B20").Value = 14475
$ws.Range("C20").Value = 125
$ws.Range("E20").Value = 5079

# Row 92 / 93: swap Ghana and Costa Rica, with Costa Rica now carrying
# updated totals and Ghana keeping its previous totals, but the two rows
# swap order in the shared list (Costa Rica now precedes Ghana).
$ws.Range("A92").Value = "Costa Rica"
$ws.Range("B92").Value = 642
$ws.Range("C92").Value = 16
$ws.Range("D92").Value = 67
$ws.Range("E92").Value = 571
$ws.Range("F92").Value = 11
$ws.Range("H92").Value = 4

$ws.Range("A93").Value = "Ghana"
$ws.Range("B93").Value = 641
$ws.Range("D93").Value = 83
$ws.Range("E93").Value = 550
$ws.Range("F93").Value = 2
$ws.Range("H93").Value = 8
